$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append "*" to all header cells in row 12 except D12 (Dependiente de quien?)
$headerCols = @("A","B","C","E","F","G","H","I","J","K","L","M","N","O","P")
foreach ($col in $headerCols) {
    $cell = $ws.Range($col + "12")
    $cell.Value = $cell.Value2 + "*"
}
$ws.Range("Q12").Value = "SE REALIZO LA LIMPIEZA* SI/NO"

# Add a comment to D12
$d12 = $ws.Range("D12")
if ($d12.Comment -ne $null) {
    $d12.Comment.Delete()
}
$comment = $d12.AddComment("RED BUCAL:`nIngrese la cedula del titular al que depende, si no depende, deje el campo vacio")

# Update selection
$ws.Range("A13:Q16").Select()
